# Fixed bug with non-greedy regular expression matching.
# Recomputed text-pattern frequency percentages in row 2 of the
# "issue_title_refactoring_doc_tex" sheet, and brought that sheet to the
# front (it was the sheet being worked on when the fix was saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("issue_title_refactoring_doc_tex")

# --- Updated data (A2:AI2) ------------------------------------------------
$ws.Range("A2").Value  = [double]"3.7792894935751997E-2"
$ws.Range("B2").Value  = [double]"0.11337868480725601"
$ws.Range("C2").Value  = [double]"0"
$ws.Range("D2").Value  = [double]"3.09901738473167"
$ws.Range("E2").Value  = [double]"0.15117157974300799"
$ws.Range("F2").Value  = [double]"29.402872260015101"
$ws.Range("G2").Value  = [double]"5.55555555555555"
$ws.Range("H2").Value  = [double]"8.6545729402872205"
$ws.Range("I2").Value  = [double]"2.87226001511715"
$ws.Range("J2").Value  = [double]"2.9856386999244098"
$ws.Range("K2").Value  = [double]"7.2184429327286397"
$ws.Range("L2").Value  = [double]"9.4860166288737702"
$ws.Range("M2").Value  = [double]"0.45351473922902402"
$ws.Range("N2").Value  = [double]"1.1715797430083099"
$ws.Range("O2").Value  = [double]"2.45653817082388"
$ws.Range("P2").Value  = [double]"3.7414965986394502"
$ws.Range("Q2").Value  = [double]"3.7792894935751997E-2"
$ws.Range("R2").Value  = [double]"0.22675736961451201"
$ws.Range("S2").Value  = [double]"1.9274376417233501"
$ws.Range("T2").Value  = [double]"0.45351473922902402"
$ws.Range("U2").Value  = [double]"9.9773242630385397"
$ws.Range("V2").Value  = [double]"1.9652305366591001"
$ws.Range("W2").Value  = [double]"0.11337868480725601"
$ws.Range("X2").Value  = [double]"3.4391534391534302"
$ws.Range("Y2").Value  = [double]"1.5495086923658301"
$ws.Range("Z2").Value  = [double]"0"
$ws.Range("AA2").Value = [double]"0.15117157974300799"
$ws.Range("AB2").Value = [double]"1.7006802721088401"
$ws.Range("AC2").Value = [double]"0.22675736961451201"
$ws.Range("AD2").Value = [double]"3.7792894935751997E-2"
$ws.Range("AE2").Value = [double]"0.11337868480725601"
$ws.Range("AF2").Value = [double]"3.7792894935751997E-2"
$ws.Range("AG2").Value = [double]"0.11337868480725601"
$ws.Range("AH2").Value = [double]"0.45351473922902402"
$ws.Range("AI2").Value = [double]"7.5585789871504105E-2"

# --- View state: bring this sheet to the front, reset the scroll position,
#     rezoom, and select the full data range as the last action before
#     saving. -----------------------------------------------------------
[void]$ws.Activate()
$excel.ActiveWindow.Zoom = 43
[void]$ws.Range("A1:AI2").Select()
